$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.657.30"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.532.78"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.34"
$ws.Range("E5").Value = "  +3.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.05"
$ws.Range("E6").Value = "  -5.95%  "
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  -3.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.21"
$ws.Range("E10").Value = "  -3.74%  "
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.921.03"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.527.99"
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.654.84"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.00"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.57"
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.13"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.21"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.55"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.21"
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.05"
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.85"
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.24"
$ws.Range("E33").Value = "  +4.06%  "
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.29"
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0785"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.111"
$ws.Range("E38").Value = "  -4.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.119"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.58"
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.30"
$ws.Range("E41").Value = "  +10.48%  "
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("E45").Value = "  -5.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.015.63"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.56"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.774.66"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.30"
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("E51").Value = "  -0.97%  "
